$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (B text, C value)
# Note: B text values start with a literal apostrophe. Excel's .Value setter
# treats a single *leading* apostrophe as a "force text" prefix marker and
# strips it, so we double it up ('' -> literal ') to get one real leading
# apostrophe written into the cell content, matching the target text exactly.
$updates = @(
    @{ Row = 2;  B = "''Akkermansia_muciniphila_ATCC_BAA_835.mat'";            C = 0 }
    @{ Row = 3;  B = "''Alistipes_finegoldii_DSM_17242.mat'";                  C = 0 }
    @{ Row = 4;  B = "''Alistipes_putredinis_DSM_17216.mat'";                  C = 0 }
    @{ Row = 5;  B = "''Bacteroides_cellulosilyticus_DSM_14838.mat'";          C = 0 }
    @{ Row = 6;  B = "''Bacteroides_coprophilus_DSM_18228.mat'";               C = 0 }
    @{ Row = 7;  B = "''Bacteroides_fragilis_3_1_12.mat'";                     C = 0 }
    @{ Row = 8;  B = "''Bacteroides_oleiciplenus_YIT_12058.mat'";              C = 0 }
    @{ Row = 9;  B = "''Bacteroides_ovatus_ATCC_8483.mat'";                    C = 0 }
    @{ Row = 10; B = "''Bacteroides_plebeius_M12_DSM_17135.mat'";              C = 0 }
    @{ Row = 11; B = "''Bacteroides_salyersiae_WAL_10018.mat'";                C = 0 }
    @{ Row = 12; B = "''Bacteroides_thetaiotaomicron_VPI_5482.mat'";           C = 0 }
    @{ Row = 13; B = "''Bacteroides_uniformis_ATCC_8492.mat'";                 C = 0 }
    @{ Row = 14; B = "''Bacteroides_vulgatus_ATCC_8482.mat'";                  C = 0 }
    @{ Row = 15; B = "''Barnesiella_intestinihominis_YIT_11860.mat'";          C = 0 }
    @{ Row = 16; B = "''Bifidobacterium_animalis_lactis_AD011.mat'";           C = 0 }
    @{ Row = 17; B = "''Bilophila_wadsworthia_3_1_6.mat'";                     C = 0 }
    @{ Row = 18; B = "''Escherichia_coli_O157_H7_str_Sakai.mat'";              C = 0 }
    @{ Row = 19; B = "''Eubacterium_limosum_KIST612.mat'";                     C = 0 }
    @{ Row = 20; B = "''Eubacterium_ramulus_ATCC_29099.mat'";                  C = 0 }
    @{ Row = 21; B = "''Flavonifractor_plautii_ATCC_29863.mat'";               C = 0.004 }
    @{ Row = 22; B = "''Marvinbryantia_formatexigens_I_52_DSM_14469.mat'";     C = 0 }
    @{ Row = 23; B = "''Odoribacter_splanchnicus_1651_6_DSM_20712.mat'";       C = 0 }
    @{ Row = 24; B = "''Parabacteroides_distasonis_ATCC_8503.mat'";            C = 0 }
    @{ Row = 25; B = "''Parabacteroides_johnsonii_DSM_18315.mat'";             C = 0 }
    @{ Row = 26; B = "''Paraprevotella_xylaniphila_YIT_11841.mat'";            C = 0 }
    @{ Row = 27; B = "''Parasutterella_excrementihominis_YIT_11859.mat'";      C = 0.032 }
    @{ Row = 28; B = "''Phascolarctobacterium_succinatutens_YIT_12067.mat'";   C = 0.075 }
    @{ Row = 29; B = "''Prevotella_copri_CB7_DSM_18205.mat'";                  C = 0 }
    @{ Row = 30; B = "''Prevotella_stercorea_DSM_18206.mat'";                  C = 0.023 }
    @{ Row = 31; B = "''Roseburia_inulinivorans_DSM_16841.mat'";               C = 0.133 }
    @{ Row = 32; B = "''Sutterella_wadsworthensis_3_1_45B.mat'";               C = 0.732 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 3).Value = $u.C
}
